$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same "想去人数" (F column) figures
# for these rows and need to be bumped to reflect the newly generated output.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2249
    $ws.Range("F4").Value = 1654
    $ws.Range("F5").Value = 7650
}
